# Apply the changes described by the commit: header columns E1:I1 get
# reordered (R_mean, R_std, RMSE_mean, RMSE_std, num_of_data) and the newly
# computed sklearn-model statistics are filled in for rows 2-9 and 83-84.
# Also updates the active selection on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("opv_results")

# --- Reorder the header row (E1:I1) ---
$ws.Range("E1").Value = "R_mean"
$ws.Range("F1").Value = "R_std"
$ws.Range("G1").Value = "RMSE_mean"
$ws.Range("H1").Value = "RMSE_std"
$ws.Range("I1").Value = "num_of_data"

# --- Fill in newly-produced model results for rows 2-9 ---
$ws.Range("E2").Value = [double]"0.62813610018402299"
$ws.Range("F2").Value = [double]"7.2833427240581297E-2"
$ws.Range("G2").Value = [double]"0.157762285934057"
$ws.Range("H2").Value = [double]"2.4092237265460802E-2"
$ws.Range("I2").Value = 556

$ws.Range("E3").Value = [double]"0.63877089008832999"
$ws.Range("F3").Value = [double]"7.9113294455477698E-2"
$ws.Range("G3").Value = [double]"0.156107845269398"
$ws.Range("H3").Value = [double]"2.55185545309627E-2"
$ws.Range("I3").Value = 556

$ws.Range("E4").Value = [double]"0.64649318977724501"
$ws.Range("F4").Value = [double]"6.0121236629882301E-2"
$ws.Range("G4").Value = [double]"0.154401391889916"
$ws.Range("H4").Value = [double]"2.3177355135618399E-2"
$ws.Range("I4").Value = 556

$ws.Range("E5").Value = [double]"0.49359585421651297"
$ws.Range("F5").Value = [double]"7.37843769701262E-2"
$ws.Range("G5").Value = [double]"0.176627652838907"
$ws.Range("H5").Value = [double]"2.2985116466043499E-2"
$ws.Range("I5").Value = 556

$ws.Range("E6").Value = [double]"0.68288723323085099"
$ws.Range("F6").Value = [double]"6.7555834775905793E-2"
$ws.Range("G6").Value = [double]"0.14637437701699699"
$ws.Range("H6").Value = [double]"1.9779066086000101E-2"
$ws.Range("I6").Value = 556

$ws.Range("E7").Value = [double]"0.675303582386243"
$ws.Range("F7").Value = [double]"4.8045653484017799E-2"
$ws.Range("G7").Value = [double]"0.14817191299363999"
$ws.Range("H7").Value = [double]"1.3441455174449699E-2"
$ws.Range("I7").Value = 556

$ws.Range("E8").Value = [double]"0.65101278918567995"
$ws.Range("F8").Value = [double]"5.6228455623029E-2"
$ws.Range("G8").Value = [double]"0.153485315711967"
$ws.Range("H8").Value = [double]"1.9843083811824199E-2"
$ws.Range("I8").Value = 556

$ws.Range("E9").Value = [double]"0.72698273872294505"
$ws.Range("F9").Value = [double]"5.5640040279394003E-2"
$ws.Range("G9").Value = [double]"0.13842668840805999"
$ws.Range("H9").Value = [double]"1.6471351077796299E-2"
$ws.Range("I9").Value = 556

# --- Fill in newly-produced model results for rows 83-84 ---
$ws.Range("E83").Value = [double]"0.55302783911153297"
$ws.Range("F83").Value = [double]"3.0489556892871798E-2"
$ws.Range("G83").Value = [double]"0.16844429999999999"
$ws.Range("H83").Value = [double]"1.2271737E-2"
$ws.Range("I83").Value = 447

$ws.Range("E84").Value = [double]"0.54921497712617395"
$ws.Range("F84").Value = [double]"3.5687978617597801E-2"
$ws.Range("G84").Value = [double]"0.16964984164206001"
$ws.Range("H84").Value = [double]"1.35076548325601E-2"
$ws.Range("I84").Value = 447

# --- Update the active selection (also clears the stale topLeftCell scroll
#     position from the previous session) ---
$ws.Activate()
$ws.Range("I17").Select()
